# Working on the goal page. Now selecting the goals in bank.
#
# This script:
#  1. Appends 9 new rows (712-720) of goal-bank entries to the "SkillBank"
#     worksheet (bringing its data range from A1:E711 to A1:E720), adding
#     7 brand-new shared strings along the way and re-using the two that
#     already existed ("Networking", "Negotiate").
#  2. Updates the selection / scroll state on a few sheets and switches
#     the active sheet from "Goal" to "SkillBank".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Add the new SkillBank rows
# ---------------------------------------------------------------------
$skillBank = $wb.Worksheets.Item("SkillBank")

# Each entry is: type id, name
$newGoals = @(
    @(6,  "Lose Weight"),
    @(6,  "Eat Healthy"),
    @(8,  "Networking"),
    @(8,  "Tracking Acomplishments "),
    @(8,  "Update your Skills"),
    @(8,  "Negotiate"),
    @(10, "Speak English fluently"),
    @(7,  "1000 subscribers to my blogs"),
    @(12, "Donate at an organisation")
)

$startRow = 712
$row = $startRow
foreach ($goal in $newGoals) {
    $id = $row - 1
    $skillBank.Cells.Item($row, 1).Value = $id
    $skillBank.Cells.Item($row, 2).Value = $goal[0]
    $skillBank.Cells.Item($row, 3).Value = $goal[1]
    $skillBank.Cells.Item($row, 4).Value = "\N"
    $skillBank.Cells.Item($row, 5).Value = "\N"
    $row++
}

# ---------------------------------------------------------------------
# 2. Update sheet view / selection state
# ---------------------------------------------------------------------

# "Goal" sheet loses the active tab, scrolls a bit further and the
# selection moves from C20 to C29.
$goalSheet = $wb.Worksheets.Item("Goal")
$goalSheet.Activate() | Out-Null
$goalSheet.Range("C29").Select() | Out-Null

# "GoalType" sheet's selection grows from the single cell C6 to the
# range A7:D14.
$goalTypeSheet = $wb.Worksheets.Item("GoalType")
$goalTypeSheet.Activate() | Out-Null
$goalTypeSheet.Range("A7:D14").Select() | Out-Null

# "SkillBank" becomes the active sheet (the tab the author is now
# working in), scrolled down near the newly-added rows, with D725
# selected.
$skillBank.Activate() | Out-Null
$skillBank.Range("D725").Select() | Out-Null
